$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

$ws.Range("D2").Value = 0.0002229288220405579
$ws.Range("E2").Value = 0.04792241798713803
$ws.Range("G2").Value = 0.004238322842866182
$ws.Range("H2").Value = 0.008408657740801573
$ws.Range("I2").Value = 0.01438153302296996
$ws.Range("J2").Value = 0.01559349102899432
$ws.Range("K2").Value = 0.00139350164681673
$ws.Range("D3").Value = 0.002423597499728203
$ws.Range("E3").Value = 0.04775964701548219
$ws.Range("G3").Value = 0.003841289319097996
$ws.Range("H3").Value = 0.01006700610741973
$ws.Range("I3").Value = 0.0127634909003973
$ws.Range("J3").Value = 0.01622256869450212
$ws.Range("K3").Value = 0.001281950157135725
$ws.Range("D4").Value = 0.003336433786898851
$ws.Range("E4").Value = 0.0581201552413404
$ws.Range("G4").Value = 0.004698055796325207
$ws.Range("H4").Value = 0.01257986901327968
$ws.Range("I4").Value = 0.01560318190604448
$ws.Range("J4").Value = 0.01919363299384713
$ws.Range("K4").Value = 0.001591862645000219
$ws.Range("D5").Value = 0.0002646297216415405
$ws.Range("E5").Value = 0.04319082386791706
$ws.Range("G5").Value = 0.003749916329979897
$ws.Range("H5").Value = 0.007721805013716221
$ws.Range("I5").Value = 0.01352619146928191
$ws.Range("J5").Value = 0.01351108774542809
$ws.Range("K5").Value = 0.001244128681719303
$ws.Range("D6").Value = 0.005630886182188988
$ws.Range("E6").Value = 0.5136187081225216
$ws.Range("G6").Value = 0.01312649948522449
$ws.Range("H6").Value = 0.03437632508575916
$ws.Range("I6").Value = 0.4128670734353364
$ws.Range("J6").Value = 0.03484869189560413
$ws.Range("K6").Value = 0.005245338659733534
$ws.Range("D8").Value = 0.0002229288220405579
$ws.Range("E8").Value = 0.04792241798713803
$ws.Range("G8").Value = 0.004238322842866182
$ws.Range("H8").Value = 0.008408657740801573
$ws.Range("I8").Value = 0.01438153302296996
$ws.Range("J8").Value = 0.01559349102899432
$ws.Range("K8").Value = 0.00139350164681673
$ws.Range("D9").Value = 0.002423597499728203
$ws.Range("E9").Value = 0.04775964701548219
$ws.Range("G9").Value = 0.003841289319097996
$ws.Range("H9").Value = 0.01006700610741973
$ws.Range("I9").Value = 0.0127634909003973
$ws.Range("J9").Value = 0.01622256869450212
$ws.Range("K9").Value = 0.001281950157135725
$ws.Range("D10").Value = 0.003336433786898851
$ws.Range("E10").Value = 0.0581201552413404
$ws.Range("G10").Value = 0.004698055796325207
$ws.Range("H10").Value = 0.01257986901327968
$ws.Range("I10").Value = 0.01560318190604448
$ws.Range("J10").Value = 0.01919363299384713
$ws.Range("K10").Value = 0.001591862645000219
$ws.Range("D11").Value = 0.0002646297216415405
$ws.Range("E11").Value = 0.04319082386791706
$ws.Range("G11").Value = 0.003749916329979897
$ws.Range("H11").Value = 0.007721805013716221
$ws.Range("I11").Value = 0.01352619146928191
$ws.Range("J11").Value = 0.01351108774542809
$ws.Range("K11").Value = 0.001244128681719303
$ws.Range("D12").Value = 0.005630886182188988
$ws.Range("E12").Value = 0.5136187081225216
$ws.Range("G12").Value = 0.01312649948522449
$ws.Range("H12").Value = 0.03437632508575916
$ws.Range("I12").Value = 0.4128670734353364
$ws.Range("J12").Value = 0.03484869189560413
$ws.Range("K12").Value = 0.005245338659733534
